# Auto-generated update of cryptos table rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to be treated as text,
# matching the original inline-string cell typing so numeric-looking
# values like "1.00" are not coerced into plain numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '69.542.79'
$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").Value = '3.769.22'
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '606.27'
$ws.Range("E5").Value = '  +1.56%  '

$ws.Range("D6").Value = '170.10'
$ws.Range("E6").Value = '  +2.54%  '

$ws.Range("D7").Value = '3.767.43'
$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +2.18%  '

$ws.Range("E10").Value = '  +4.41%  '

$ws.Range("E11").Value = '  +3.92%  '

$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").Value = '38.65'
$ws.Range("E13").Value = '  +2.54%  '

$ws.Range("E14").Value = '  +4.42%  '

$ws.Range("D15").Value = '4.402.40'
$ws.Range("E15").Value = '  +1.69%  '

$ws.Range("D16").Value = '3.768.59'
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").Value = '69.504.93'
$ws.Range("E17").Value = '  +2.95%  '

$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").Value = '17.19'
$ws.Range("E20").Value = '  -1.73%  '

$ws.Range("D21").Value = '10.92'
$ws.Range("E21").Value = '  +16.62%  '

$ws.Range("D22").Value = '497.79'
$ws.Range("E22").Value = '  +1.22%  '

$ws.Range("E23").Value = '  +1.28%  '

$ws.Range("E24").Value = '  +11.73%  '

$ws.Range("D25").Value = '85.68'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("E26").Value = '  +2.36%  '

$ws.Range("D27").Value = '12.45'
$ws.Range("E27").Value = '  +2.14%  '

$ws.Range("E28").Value = '  +2.50%  '

$ws.Range("E29").Value = '  -0.19%  '

$ws.Range("E30").Value = '  +2.49%  '

$ws.Range("E31").Value = '  +6.92%  '

$ws.Range("D32").Value = '8.06'
$ws.Range("E32").Value = '  +5.70%  '

$ws.Range("D33").Value = '32.12'
$ws.Range("E33").Value = '  +1.83%  '

$ws.Range("D34").Value = '3.915.55'
$ws.Range("E34").Value = '  +1.81%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.109'
$ws.Range("E35").Value = '  +1.55%  '

$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.703.55'
$ws.Range("E36").Value = '  +1.40%  '

$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.09%  '

$ws.Range("E38").Value = '  +1.98%  '

$ws.Range("E39").Value = '  +3.26%  '

$ws.Range("E40").Value = '  +2.29%  '

$ws.Range("E41").Value = '  +1.57%  '

$ws.Range("D42").Value = '3.09'
$ws.Range("E42").Value = '  +10.98%  '

$ws.Range("D43").Value = '442.75'
$ws.Range("E43").Value = '  +2.29%  '

$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '48.75'
$ws.Range("E44").Value = '  +0.28%  '

$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").Value = '2.00'
$ws.Range("E45").Value = '  +3.30%  '

$ws.Range("E46").Value = '  +1.45%  '

$ws.Range("D48").Value = '40.81'
$ws.Range("E48").Value = '  +0.21%  '

$ws.Range("D49").Value = '2.817.74'
$ws.Range("E49").Value = '  +2.24%  '

$ws.Range("D50").Value = '140.84'
$ws.Range("E50").Value = '  -1.36%  '

$ws.Range("E51").Value = '  +2.66%  '
